$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: the "varEst" parameter's Type changes from "est" to "prep"
$ws.Cells.Item(25, 3).Value = "prep"

# New rows 26-27: north vector X/Y coordinates.
# Variable names first, then labels, to match shared-string append order.
$ws.Cells.Item(26, 2).Value = "northx"
$ws.Cells.Item(27, 2).Value = "northy"
$ws.Cells.Item(26, 1).Value = "X cooridante for north vector"
$ws.Cells.Item(27, 1).Value = "Y cooridante for north vector"

$ws.Cells.Item(26, 3).Value = "data"
$ws.Cells.Item(26, 4).Value = 0
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 0

$ws.Cells.Item(27, 3).Value = "data"
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(27, 5).Value = 1
$ws.Cells.Item(27, 6).Value = 1

# Update the selected/active cell to match the authored state
$ws.Range("B17").Select()
